$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row's formatting (borders / wrap / style index) down to the
# two new rows before filling in values, so every new cell carries the same
# bordered "s=2" style the header row uses.
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K3").PasteSpecial(-4122)

# Row 2 - Botswana AMR Gram Stain Observation
$ws.Range("A2").Value = "botswana-amr-gram-stain-observation"
$ws.Range("B2").Value = "Botswana AMR Gram Stain Observation"
$ws.Range("E2").Value = "LOINC#664-3"
$ws.Range("G2").Value = "dateTime"
$ws.Range("H2").Value = "CodeableConcept"
$ws.Range("I2").Value = "optional"

# Row 3 - Botswana AMR Organism Identification Observation
$ws.Range("A3").Value = "botswana-amr-organism-observation"
$ws.Range("B3").Value = "Botswana AMR Organism Identification Observation"
$ws.Range("E3").Value = "LOINC#634-6"
$ws.Range("G3").Value = "dateTime"
$ws.Range("H3").Value = "CodeableConcept"
$ws.Range("I3").Value = "optional"
